# pptx: Fix list level numbering
#
# All paragraphs that were (incorrectly) written one level too deep
# (PowerPoint COM IndentLevel == 2, i.e. OOXML a:pPr lvl="1") are
# promoted back to the top level (IndentLevel == 1, i.e. lvl="0").
#
# This mirrors the fix described in the commit: a top-level list item
# should be at the same outline level as a top-level paragraph - only
# continuation levels of a list should be indented further.

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)

    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)

        if ($shape.HasTextFrame) {
            $tf = $shape.TextFrame
            $tr = $tf.TextRange
            $count = $tr.Paragraphs().Count

            for ($pi = 1; $pi -le $count; $pi++) {
                $para = $tr.Paragraphs($pi, 1)
                if ($para.IndentLevel -eq 2) {
                    $para.IndentLevel = 1
                }
            }
        }
    }
}
